$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Price-column values are plain text in the sheet (e.g. "63.639.12",
# "568.98"), so a leading apostrophe is used to force text entry and
# keep the original "General" cell format/type instead of letting Excel
# reinterpret them as numbers.

$ws.Range("D2").Value = "'63.639.12"
$ws.Range("E2").Value = '  -1.50%  '

$ws.Range("D3").Value = "'3.406.71"
$ws.Range("E3").Value = '  -0.43%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'568.98"
$ws.Range("E5").Value = '  -0.70%  '

$ws.Range("D6").Value = "'157.41"
$ws.Range("E6").Value = '  +0.34%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = "'3.409.03"
$ws.Range("E8").Value = '  -0.50%  '

$ws.Range("D9").Value = "'0.565"
$ws.Range("E9").Value = '  -10.05%  '

$ws.Range("E10").Value = '  +1.01%  '

$ws.Range("E11").Value = '  -3.96%  '

$ws.Range("D12").Value = "'0.421"

$ws.Range("D13").Value = "'3.993.53"
$ws.Range("E13").Value = '  -0.53%  '

$ws.Range("E14").Value = '  +0.01%  '

$ws.Range("D15").Value = "'27.03"
$ws.Range("E15").Value = '  -3.03%  '

$ws.Range("E16").Value = '  -8.61%  '

$ws.Range("D17").Value = "'63.708.38"
$ws.Range("E17").Value = '  -1.37%  '

$ws.Range("D18").Value = "'3.422.53"
$ws.Range("E18").Value = '  -0.86%  '

$ws.Range("D19").Value = "'6.08"
$ws.Range("E19").Value = '  -4.22%  '

$ws.Range("E20").Value = '  -2.66%  '

$ws.Range("D21").Value = "'383.57"
$ws.Range("E21").Value = '  +1.47%  '

$ws.Range("E22").Value = '  -3.43%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("D24").Value = "'71.14"
$ws.Range("E24").Value = '  -1.83%  '

$ws.Range("E25").Value = '  -6.41%  '

$ws.Range("E26").Value = '  -4.58%  '

$ws.Range("D27").Value = "'9.66"
$ws.Range("E27").Value = '  -6.52%  '

$ws.Range("E28").Value = '  +0.16%  '

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  -0.08%  '

$ws.Range("E30").Value = '  -1.77%  '

$ws.Range("E31").Value = '  -6.78%  '

$ws.Range("E32").Value = '  -2.26%  '

$ws.Range("E33").Value = '  +0.03%  '

$ws.Range("E34").Value = '  -1.16%  '

$ws.Range("D35").Value = "'6.93"
$ws.Range("E35").Value = '  -3.96%  '

$ws.Range("E36").Value = '  -6.53%  '

$ws.Range("E37").Value = '  +0.67%  '

$ws.Range("E38").Value = '  +9.20%  '

$ws.Range("E39").Value = '  -4.58%  '

$ws.Range("D40").Value = "'2.792.80"
$ws.Range("E40").Value = '  -3.22%  '

$ws.Range("E41").Value = '  -5.39%  '

$ws.Range("D42").Value = "'25.85"
$ws.Range("E42").Value = '  -3.73%  '

$ws.Range("E43").Value = '  +0.11%  '

$ws.Range("D44").Value = "'25.94"
$ws.Range("E44").Value = '  -2.78%  '

$ws.Range("D45").Value = "'6.36"
$ws.Range("E45").Value = '  -8.61%  '

$ws.Range("E46").Value = '  -6.56%  '

$ws.Range("E47").Value = '  -4.97%  '

$ws.Range("D48").Value = "'2.36"
$ws.Range("E48").Value = '  +8.20%  '

$ws.Range("D49").Value = "'326.76"
$ws.Range("E49").Value = '  +1.48%  '

$ws.Range("E50").Value = '  -5.27%  '

$ws.Range("E51").Value = '  -5.78%  '
